# "Added full transfer experiments" - fills in the previously-empty
# transfer-learning result cells (columns C..K) for rows 17 and 18,
# matching the author's original entry order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17 ---
$ws.Range("J17").Value = "58.8 (0.0)"
$ws.Range("E17").Value = "36.6 (4.6)"
$ws.Range("G17").Value = "X12012310"
$ws.Range("K17").Value = "4.6 (0.4)"
$ws.Range("F17").Value = "5.2 (0.5)"
$ws.Range("D17").Value = "28.3 (6.7)"
$ws.Range("H17").Value = "41.7 (7.7)"
$ws.Range("C17").Value = "52.3 (1.1)"
$ws.Range("I17").Value = "100.0 (0)"

# --- Row 18 ---
$ws.Range("J18").Value = "19.4 (0.8)"
$ws.Range("E18").Value = "19.4 (2.4)"
$ws.Range("F18").Value = "32.4 (0.4)"
$ws.Range("K18").Value = "28.0 (0.8)"
$ws.Range("D18").Value = "48.1 (18.6)"
$ws.Range("I18").Value = "77.8 (3.7)"
$ws.Range("C18").Value = "12.0 (0.7)"
$ws.Range("H18").Value = "11.3 (0.3)"
$ws.Range("G18").Value = "X11002310"

# Restore the saved selection state
$ws.Range("H18").Select()
